$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: E31 set contents reordered, F31 Neutral -> Win
$ws.Range("E31").Value2 = "{'str', 'any'}"
$ws.Range("F31").Value2 = "Win"
$ws.Range("F31").Interior.Color = 32768

# Row 32: E32 any -> str, F32 Neutral -> Win
$ws.Range("E32").Value2 = "str"
$ws.Range("F32").Value2 = "Win"
$ws.Range("F32").Interior.Color = 32768

# Row 92: E92 set contents reordered, F92 Neutral -> Win
$ws.Range("E92").Value2 = "{'Tuple[str]', 'str', 'any'}"
$ws.Range("F92").Value2 = "Win"
$ws.Range("F92").Interior.Color = 32768

# Row 93: E93 any -> Tuple[str], F93 Neutral -> Win
$ws.Range("E93").Value2 = "Tuple[str]"
$ws.Range("F93").Value2 = "Win"
$ws.Range("F93").Interior.Color = 32768

# Row 96: F96 Neutral -> Win (E96 text unchanged)
$ws.Range("F96").Value2 = "Win"
$ws.Range("F96").Interior.Color = 32768

# Row 165: Scalpel Wins count 8 -> 13
$ws.Range("F165").Value2 = 13

# Row 166: add new "Scalpel Accuracy:" label + value in C166/D166,
# and clear E166/F166 since the "Accuracy over PyType" summary that
# used to live there moves down to row 167 (A166/B166 stay blank)
$ws.Range("C166").Value2 = "Scalpel Accuracy:"
$ws.Range("D166").Value2 = 2616.67
$ws.Range("E166").Value2 = ""
$ws.Range("F166").Value2 = ""

# Row 167 (new row): the "Accuracy over PyType" summary moves here,
# with its value updated; touch A167:D167 so they exist as blank
# styled cells matching the row above
$ws.Range("A167:F167").Interior.Color = 16777215
$ws.Range("E167").Value2 = "Accuracy over PyType"
$ws.Range("F167").Value2 = 216.67
